$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 463.732605
$ws.Range("H2").Value = 1391.197815
$ws.Range("I2").Value = 0.3632113435366598
$ws.Range("J2").Value = 0.3632113435366598
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 7.248785666666667
$ws.Range("N2").Value = 21.746357
$ws.Range("O2").Value = 0.07891374419744837
$ws.Range("P2").Value = 0.07891374419744837
$ws.Range("Q2").Value = 3361.498260289995
$ws.Range("R2").Value = 30253.48434260996
$ws.Range("S2").Value = 0.02866236705346352
$ws.Range("T2").Value = 0.02866236705346352
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 463.732605
$ws.Range("H3").Value = 1391.197815
$ws.Range("I3").Value = 0.3632113435366598
$ws.Range("J3").Value = 0.3632113435366598
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 37.82684066666667
$ws.Range("N3").Value = 113.480522
$ws.Range("O3").Value = 0.4118010609547572
$ws.Range("P3").Value = 0.4118010609547572
$ws.Range("Q3").Value = 17541.53936127327
$ws.Range("R3").Value = 157873.8542514594
$ws.Range("S3").Value = 0.1495708166191993
$ws.Range("T3").Value = 0.1495708166191993
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 463.732605
$ws.Range("H4").Value = 1391.197815
$ws.Range("I4").Value = 0.3632113435366598
$ws.Range("J4").Value = 0.3632113435366598
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 38.20927633333334
$ws.Range("N4").Value = 114.627829
$ws.Range("O4").Value = 0.4159644383477588
$ws.Range("P4").Value = 0.4159644383477588
$ws.Range("Q4").Value = 17718.88724922152
$ws.Range("R4").Value = 159469.9852429936
$ws.Range("S4").Value = 0.1510830025157616
$ws.Range("T4").Value = 0.1510830025157616
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 463.732605
$ws.Range("H5").Value = 1391.197815
$ws.Range("I5").Value = 0.3632113435366598
$ws.Range("J5").Value = 0.3632113435366598
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.572171666666666
$ws.Range("N5").Value = 25.716515
$ws.Range("O5").Value = 0.09332075650003555
$ws.Range("P5").Value = 0.09332075650003555
$ws.Range("Q5").Value = 3975.195497490524
$ws.Range("R5").Value = 35776.75947741472
$ws.Range("S5").Value = 0.03389515734823539
$ws.Range("T5").Value = 0.03389515734823539
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 12.24662533333333
$ws.Range("H6").Value = 36.739876
$ws.Range("I6").Value = 0.009591978638444229
$ws.Range("J6").Value = 0.009591978638444227
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.248785666666667
$ws.Range("N6").Value = 21.746357
$ws.Range("O6").Value = 0.07891374419744837
$ws.Range("P6").Value = 0.07891374419744837
$ws.Range("Q6").Value = 88.77316218130358
$ws.Range("R6").Value = 798.9584596317321
$ws.Range("S6").Value = 0.000756938948621577
$ws.Range("T6").Value = 0.0007569389486215769
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 12.24662533333333
$ws.Range("H7").Value = 36.739876
$ws.Range("I7").Value = 0.009591978638444229
$ws.Range("J7").Value = 0.009591978638444227
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 37.82684066666667
$ws.Range("N7").Value = 113.480522
$ws.Range("O7").Value = 0.4118010609547572
$ws.Range("P7").Value = 0.4118010609547572
$ws.Range("Q7").Value = 463.2511451883636
$ws.Range("R7").Value = 4169.260306695272
$ws.Range("S7").Value = 0.003949986979966701
$ws.Range("T7").Value = 0.003949986979966701
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 12.24662533333333
$ws.Range("H8").Value = 36.739876
$ws.Range("I8").Value = 0.009591978638444229
$ws.Range("J8").Value = 0.009591978638444227
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 38.20927633333334
$ws.Range("N8").Value = 114.627829
$ws.Range("O8").Value = 0.4159644383477588
$ws.Range("P8").Value = 0.4159644383477588
$ws.Range("Q8").Value = 467.9346915121338
$ws.Range("R8").Value = 4211.412223609204
$ws.Range("S8").Value = 0.003989922006984154
$ws.Range("T8").Value = 0.003989922006984154
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 12.24662533333333
$ws.Range("H9").Value = 36.739876
$ws.Range("I9").Value = 0.009591978638444229
$ws.Range("J9").Value = 0.009591978638444227
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 8.572171666666666
$ws.Range("N9").Value = 25.716515
$ws.Range("O9").Value = 0.09332075650003555
$ws.Range("P9").Value = 0.09332075650003555
$ws.Range("Q9").Value = 104.9801746946822
$ws.Range("R9").Value = 944.82157225214
$ws.Range("S9").Value = 0.0008951307028717964
$ws.Range("T9").Value = 0.0008951307028717962
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 689.7685036666667
$ws.Range("H10").Value = 2069.305511
$ws.Range("I10").Value = 0.5402504422695089
$ws.Range("J10").Value = 0.5402504422695089
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 7.248785666666667
$ws.Range("N10").Value = 21.746357
$ws.Range("O10").Value = 0.07891374419744837
$ws.Range("P10").Value = 0.07891374419744837
$ws.Range("Q10").Value = 4999.984042697048
$ws.Range("R10").Value = 44999.85638427344
$ws.Range("S10").Value = 0.04263318520381437
$ws.Range("T10").Value = 0.04263318520381437
# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 689.7685036666667
$ws.Range("H11").Value = 2069.305511
$ws.Range("I11").Value = 0.5402504422695089
$ws.Range("J11").Value = 0.5402504422695089
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 37.82684066666667
$ws.Range("N11").Value = 113.480522
$ws.Range("O11").Value = 0.4118010609547572
$ws.Range("P11").Value = 0.4118010609547572
$ws.Range("Q11").Value = 26091.76328508409
$ws.Range("R11").Value = 234825.8695657568
$ws.Range("S11").Value = 0.2224757053078606
$ws.Range("T11").Value = 0.2224757053078606
# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 689.7685036666667
$ws.Range("H12").Value = 2069.305511
$ws.Range("I12").Value = 0.5402504422695089
$ws.Range("J12").Value = 0.5402504422695089
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 38.20927633333334
$ws.Range("N12").Value = 114.627829
$ws.Range("O12").Value = 0.4159644383477588
$ws.Range("P12").Value = 0.4159644383477588
$ws.Range("Q12").Value = 26355.55536262952
$ws.Range("R12").Value = 237199.9982636656
$ws.Range("S12").Value = 0.2247249717857646
$ws.Range("T12").Value = 0.2247249717857646
# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 689.7685036666667
$ws.Range("H13").Value = 2069.305511
$ws.Range("I13").Value = 0.5402504422695089
$ws.Range("J13").Value = 0.5402504422695089
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 8.572171666666666
$ws.Range("N13").Value = 25.716515
$ws.Range("O13").Value = 0.09332075650003555
$ws.Range("P13").Value = 0.09332075650003555
$ws.Range("Q13").Value = 5912.814023690462
$ws.Range("R13").Value = 53215.32621321416
$ws.Range("S13").Value = 0.05041657997206935
$ws.Range("T13").Value = 0.05041657997206935
# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 111.00921
$ws.Range("H14").Value = 333.02763
$ws.Range("I14").Value = 0.08694623555538696
$ws.Range("J14").Value = 0.08694623555538696
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 7.248785666666667
$ws.Range("N14").Value = 21.746357
$ws.Range("O14").Value = 0.07891374419744837
$ws.Range("P14").Value = 0.07891374419744837
$ws.Range("Q14").Value = 804.68197031599
$ws.Range("R14").Value = 7242.137732843911
$ws.Range("S14").Value = 0.006861252991548897
$ws.Range("T14").Value = 0.006861252991548897
# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 111.00921
$ws.Range("H15").Value = 333.02763
$ws.Range("I15").Value = 0.08694623555538696
$ws.Range("J15").Value = 0.08694623555538696
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 37.82684066666667
$ws.Range("N15").Value = 113.480522
$ws.Range("O15").Value = 0.4118010609547572
$ws.Range("P15").Value = 0.4118010609547572
$ws.Range("Q15").Value = 4199.12769920254
$ws.Range("R15").Value = 37792.14929282286
$ws.Range("S15").Value = 0.03580455204773058
$ws.Range("T15").Value = 0.03580455204773058
# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 111.00921
$ws.Range("H16").Value = 333.02763
$ws.Range("I16").Value = 0.08694623555538696
$ws.Range("J16").Value = 0.08694623555538696
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 38.20927633333334
$ws.Range("N16").Value = 114.627829
$ws.Range("O16").Value = 0.4159644383477588
$ws.Range("P16").Value = 0.4159644383477588
$ws.Range("Q16").Value = 4241.58158043503
$ws.Range("R16").Value = 38174.23422391527
$ws.Range("S16").Value = 0.03616654203924848
$ws.Range("T16").Value = 0.03616654203924848
# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 111.00921
$ws.Range("H17").Value = 333.02763
$ws.Range("I17").Value = 0.08694623555538696
$ws.Range("J17").Value = 0.08694623555538696
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 8.572171666666666
$ws.Range("N17").Value = 25.716515
$ws.Range("O17").Value = 0.09332075650003555
$ws.Range("P17").Value = 0.09332075650003555
$ws.Range("Q17").Value = 951.5900047010499
$ws.Range("R17").Value = 951.5900047010499
$ws.Range("S17").Value = 0.008113888476858998
$ws.Range("T17").Value = 0.008113888476858998
